$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: 综合准班率指数(%)
$ws.Range("B2").Value = "'35.79"
$ws.Range("C2").Value = "'2024-10"

# Update row 3: 到离港服务准班率指数(%)
$ws.Range("B3").Value = "'30.55"
$ws.Range("C3").Value = "'2024-10"

# Update row 4: 收发货服务准班率指数(%)
$ws.Range("B4").Value = "'41.04"
$ws.Range("C4").Value = "'2024-10"
